$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.001754606491886079
$ws.Range("B2").Value = 0.001777959405444562
$ws.Range("C2").Value = 0.001120834378525615
$ws.Range("D2").Value = 0.004016589373350143
$ws.Range("E2").Value = 0.0008744585793465376
$ws.Range("F2").Value = 0.0002035632642218843
$ws.Range("G2").Value = 0.001417377614416182
$ws.Range("H2").Value = 0.0003286282881163061
$ws.Range("I2").Value = 0.00006886640039738268
$ws.Range("J2").Value = 0.002073989948257804
$ws.Range("K2").Value = 0.001416653976775706
$ws.Range("L2").Value = 0.0007442363421432674
$ws.Range("M2").Value = 0.004606277216225863
$ws.Range("N2").Value = 0.0007190199103206396
$ws.Range("O2").Value = 0.002032072050496936
$ws.Range("P2").Value = 0.000147491751704365
$ws.Range("Q2").Value = 0.001229967689141631
$ws.Range("R2").Value = 0.0009951511165127158
$ws.Range("S2").Value = 0.000429465901106596
$ws.Range("T2").Value = 0.0008939016843214631
$ws.Range("U2").Value = 0.0006867075571790338
$ws.Range("V2").Value = 0.0004038789193145931
$ws.Range("W2").Value = 0.001617243397049606
$ws.Range("X2").Value = 0.00101487641222775
$ws.Range("Y2").Value = 0.0005453828489407897
$ws.Range("Z2").Value = 0.001321452204138041
$ws.Range("AA2").Value = 0.0000236046253121458
$ws.Range("AB2").Value = 0.0006420247373171151
$ws.Range("AC2").Value = 0.001132911304011941
$ws.Range("AD2").Value = 0.0008127611363306642
$ws.Range("AE2").Value = 0.0004678583063650876
$ws.Range("AF2").Value = 0.0008771519642323256
$ws.Range("AG2").Value = 0.0009533445117995143
$ws.Range("AH2").Value = 0.0005029041785746813
$ws.Range("AI2").Value = 0.00006003598173265345
$ws.Range("AJ2").Value = 0.0006728211883455515
$ws.Range("AK2").Value = 0.00041783886263147
$ws.Range("AL2").Value = 0.001761141582392156
$ws.Range("AM2").Value = 0.0003719684609677643
$ws.Range("AN2").Value = 0.0008952082716859877
$ws.Range("AO2").Value = 0.0003177184844389558
$ws.Range("AP2").Value = 0.0005896556540392339
$ws.Range("AQ2").Value = 0.0001200065162265673
$ws.Range("AR2").Value = 0.00035045703407377
$ws.Range("AS2").Value = 0.0009188601397909224
$ws.Range("AT2").Value = 0.0006259909132495522
$ws.Range("AU2").Value = 0.001187238725833595
$ws.Range("AV2").Value = 0.0005939698312431574
$ws.Range("AW2").Value = 0.0006781402626074851
$ws.Range("AX2").Value = 0.0001735024125082418
$ws.Range("AY2").Value = 0.001002419972792268
$ws.Range("AZ2").Value = 0.0001573318731971085
$ws.Range("BA2").Value = 0.0001059078495018184
$ws.Range("BB2").Value = 0.0003473657998256385
$ws.Range("BC2").Value = 0.0003733294433914125
$ws.Range("BD2").Value = 0.0003456552512943745
$ws.Range("BE2").Value = 0.0002767109544947743
$ws.Range("BF2").Value = 0.001764470944181085
$ws.Range("BG2").Value = 0.0005129294586367905
$ws.Range("BH2").Value = 0.0001839228498283774
$ws.Range("BI2").Value = 0.0007951530278660357
$ws.Range("BJ2").Value = 0.000358313147444278
$ws.Range("BK2").Value = 0.001116806059144437
$ws.Range("BL2").Value = 0.001930131809785962
$ws.Range("BM2").Value = 0.001150695024989545
$ws.Range("BN2").Value = 0.00001787976361811161
$ws.Range("BO2").Value = 0.0003838514385279268
$ws.Range("BP2").Value = 0.000222015893086791
$ws.Range("BQ2").Value = 0.000619741331320256
$ws.Range("BR2").Value = 0.0003246865817345679
$ws.Range("BS2").Value = 0.0003255271585658193
$ws.Range("BT2").Value = 0.0002076693635899574
$ws.Range("BU2").Value = 0.0007859764155000448
$ws.Range("BV2").Value = 0.0006162981153465807
$ws.Range("BW2").Value = 0.000633406569249928
$ws.Range("BX2").Value = 0.00007246059976750985
$ws.Range("BY2").Value = 0.001380055211484432
$ws.Range("BZ2").Value = 0.00006326096627162769
$ws.Range("CA2").Value = 0.0002922783023677766
$ws.Range("CB2").Value = 0.0003588296822272241
$ws.Range("CC2").Value = 0.0004726642509922385
$ws.Range("CD2").Value = 0.001174705103039742
$ws.Range("CE2").Value = 0.0002622714673634619
$ws.Range("CF2").Value = 0.0001887840771814808
$ws.Range("CG2").Value = 0.0009994104038923979
$ws.Range("CH2").Value = 0.0002302034408785403
$ws.Range("CI2").Value = 0.0003668310819193721
$ws.Range("CJ2").Value = 0.0000996698800008744
$ws.Range("CK2").Value = 0.0001991850294871256
$ws.Range("CL2").Value = 0.0001050188730005175
$ws.Range("CM2").Value = 0.0003022708697244525
$ws.Range("CN2").Value = 0.00008520515984855592
$ws.Range("CO2").Value = 0.0002107712207362056
$ws.Range("CP2").Value = 0.0001199797552544624
$ws.Range("CQ2").Value = 0.000919884187169373
$ws.Range("CR2").Value = 0.0002438901283312589
$ws.Range("CS2").Value = 0.0008630980737507343
$ws.Range("CT2").Value = 0.0005573917296715081
$ws.Range("CU2").Value = 0.00130100769456476
$ws.Range("CV2").Value = 0.001203915802761912
$ws.Range("CW2").Value = 0.0001417133753420785
$ws.Range("CX2").Value = 0.00009514174598734826
$ws.Range("CY2").Value = 0.0007860083132982254
$ws.Range("CZ2").Value = 0.00002132550071110018
$ws.Range("DA2").Value = 0.0004990854067727923
$ws.Range("DB2").Value = 0.00007565418491140008
$ws.Range("DC2").Value = 0.0002855685888789594
$ws.Range("DD2").Value = 0.00001387480642733863
$ws.Range("DE2").Value = 0.0001204506115755066
$ws.Range("DF2").Value = 0.001313401036895812
$ws.Range("DG2").Value = 0.0007641489501111209
$ws.Range("DH2").Value = 0.00009271075396100059
$ws.Range("DI2").Value = 0.0002256590960314497
$ws.Range("DJ2").Value = 0.0002752221480477601
$ws.Range("DK2").Value = 0.001907362486235797
$ws.Range("DL2").Value = 0.00001194480864796788
$ws.Range("DM2").Value = 0.001914367661811411
$ws.Range("DN2").Value = 0.000683351478073746
$ws.Range("DO2").Value = 0.0004864971269853413
$ws.Range("DP2").Value = 0.001112466095946729
$ws.Range("DQ2").Value = 0.001265911851078272
$ws.Range("DR2").Value = 0.001216406933963299
$ws.Range("DS2").Value = 0.0007214384968392551
$ws.Range("DT2").Value = 0.0005101384012959898
$ws.Range("DU2").Value = 0.0001702255103737116
$ws.Range("DV2").Value = 0.0007457251776941121
$ws.Range("DW2").Value = 0.001192152965813875
$ws.Range("DX2").Value = 0.001326717901974916
$ws.Range("DY2").Value = 0.0009123856434598565
$ws.Range("DZ2").Value = 0.001124119269661605
$ws.Range("EA2").Value = 0.00004607994924299419
$ws.Range("EB2").Value = 0.001005593920126557
$ws.Range("EC2").Value = 0.000283252855297178
$ws.Range("ED2").Value = 0.0009491195087321103
$ws.Range("EE2").Value = 0.0007719402783550322
$ws.Range("EF2").Value = 0.0003655496693681926
$ws.Range("EG2").Value = 0.0001508185087004676
$ws.Range("EH2").Value = 0.0001046046818373725
$ws.Range("EI2").Value = 0.002480720169842243
$ws.Range("EJ2").Value = 0.0002462422417011112
$ws.Range("EK2").Value = 0.0002834755287040025
$ws.Range("EL2").Value = 0.00003888938226737082
$ws.Range("EM2").Value = 0.000944889266975224
$ws.Range("EN2").Value = 0.00043240882223472
$ws.Range("EO2").Value = 0.001323794946074486
$ws.Range("EP2").Value = 0.00006786978337913752
$ws.Range("EQ2").Value = 0.0007716879481449723
$ws.Range("ER2").Value = 0.0004665349260903895
$ws.Range("ES2").Value = 0.0001469024282414466
$ws.Range("ET2").Value = 0.00004193518179818057
$ws.Range("EU2").Value = 0.0002885486173909158
$ws.Range("EV2").Value = 0.001503924140706658
$ws.Range("EW2").Value = 0.0004500369832385331
$ws.Range("EX2").Value = 0.001316034467890859
$ws.Range("EY2").Value = 0.0002261569607071579
$ws.Range("EZ2").Value = 0.00006943259359104559
$ws.Range("FA2").Value = 0.001024012337438762
$ws.Range("FB2").Value = 0.00007581402314826846
$ws.Range("FC2").Value = 0.0001964737894013524
$ws.Range("FD2").Value = 0.00003541555997799151
$ws.Range("FE2").Value = 0.00008554661326343194
$ws.Range("FF2").Value = 0.00005064260403742082
$ws.Range("FG2").Value = 0.00001296796835958958
$ws.Range("FH2").Value = 0.001463005552068353
$ws.Range("FI2").Value = 0.0002945977030321956
$ws.Range("FJ2").Value = 0.0006744017591699958
$ws.Range("FK2").Value = 0.00135689543094486
$ws.Range("FL2").Value = 0.0008936421945691109
$ws.Range("FM2").Value = 0.0004771417588926852
$ws.Range("FN2").Value = 0.00004642961357603781
$ws.Range("FO2").Value = 0.00006260833470150828
$ws.Range("FP2").Value = 0.002268152544274926
$ws.Range("FQ2").Value = 0.0004994447808712721
$ws.Range("FR2").Value = 0.0002374115429120138
$ws.Range("FS2").Value = 0.0003556277661118656
$ws.Range("FT2").Value = 0.00007413119601551443
$ws.Range("FU2").Value = 0.0009620258933864534
$ws.Range("FV2").Value = 0.003669213736429811
$ws.Range("FW2").Value = 0.002440489362925291
$ws.Range("FX2").Value = 0.0002953282673843205
$ws.Range("FY2").Value = 0.001762636820785701
$ws.Range("FZ2").Value = 0.000279879430308938
$ws.Range("GA2").Value = 0.001507606822997332
$ws.Range("GB2").Value = 0.0007500199717469513
$ws.Range("GC2").Value = 0.0005088732577860355
$ws.Range("GD2").Value = 0.001615565735846758
$ws.Range("GE2").Value = 0.001981242327019572
$ws.Range("GF2").Value = 0.001520392019301653
$ws.Range("GG2").Value = 0.001355156069621444
